# 6. daļa (rediģēšanas funkcija)
$wb = $excel.ActiveWorkbook

# --- "Limiti" sheet: update the monthly-limit / remainder values ---
$limiti = $wb.Worksheets.Item("Limiti")

$newValues = @{
    2  = 200
    3  = 30
    4  = 450
    5  = 123
    6  = 242
    7  = 12
    8  = 12
    9  = 200.3
    10 = 12.2
    11 = 13.2
    12 = 1.99
    13 = 3
}

foreach ($row in $newValues.Keys) {
    $value = $newValues[$row]
    $limiti.Range("B$row").Value = $value
    $limiti.Range("C$row").Value = $value
}

# --- "Izdevumi" sheet: split the combined "KategorijaIzdevums" header into two columns ---
$izdevumi = $wb.Worksheets.Item("Izdevumi")

$izdevumi.Range("B1").Value = "Kategorija"
$izdevumi.Range("C1").Value = "Izdevums"

# Match the header formatting used by the rest of row 1 (bold, centered, bordered).
$izdevumi.Range("B1").Copy()
$izdevumi.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
